$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = 45973
$ws.Range("B26").Value = 680
$ws.Range("C26").Value = 26
$ws.Range("D26").Value = 654

$ws.Range("A26:D26").Select()
